# Auto-generated edit script applying numeric updates per the commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 3610.838
$ws.Range("J17").Value = 3854.147
$ws.Range("L17").Value = 11562.441
$ws.Range("N17").Value = -11898.441

$ws.Range("H19").Value = 3091.3157
$ws.Range("I19").Value = 1049.5555
$ws.Range("J19").Value = 4928.9
$ws.Range("K19").Value = 1049.5555
$ws.Range("L19").Value = 4928.9
$ws.Range("M19").Value = -874.5554999999999
$ws.Range("N19").Value = -5278.9

$ws.Range("H34").Value = 1309.5555
$ws.Range("I34").Value = 1309.5555
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 1309.5555
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -1106.5555
$ws.Range("N34").ClearContents()

$ws.Range("H36").Value = 1309.5555
$ws.Range("I36").Value = 1309.5555
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 1309.5555
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = -594.5554999999999
$ws.Range("N36").ClearContents()

$ws.Range("H69").Value = 500050000
$ws.Range("J69").Value = 1000000000
$ws.Range("L69").Value = 3000000000
$ws.Range("N69").Value = -3000001748

$ws.Range("H72").Value = 500050000
$ws.Range("J72").Value = 1000000000
$ws.Range("L72").Value = 9000000000
$ws.Range("N72").Value = -9000008736

$ws.Range("H86").Value = 4037.125
$ws.Range("I86").Value = 4083
$ws.Range("J86").Value = 3899.5
$ws.Range("K86").Value = 4083
$ws.Range("L86").Value = 3899.5
$ws.Range("M86").Value = -2960
$ws.Range("N86").Value = -6145.5

$ws.Range("H89").Value = 4037.125
$ws.Range("I89").Value = 4083
$ws.Range("J89").Value = 3899.5
$ws.Range("K89").Value = 20415
$ws.Range("L89").Value = 19497.5
$ws.Range("M89").Value = -14799
$ws.Range("N89").Value = -30729.5

$ws.Range("H121").Value = 3999
$ws.Range("J121").Value = 3999
$ws.Range("L121").Value = 11997
$ws.Range("N121").Value = -15491

$ws.Range("H137").Value = 68617.55499999999
$ws.Range("I137").Value = 139059.61
$ws.Range("J137").Value = 3207.0715
$ws.Range("K137").Value = 417178.83
$ws.Range("L137").Value = 9621.2145
$ws.Range("M137").Value = -414628.83
$ws.Range("N137").Value = -14721.2145

$ws.Range("H138").Value = 3311.5417
$ws.Range("J138").Value = 3461.4048
$ws.Range("L138").Value = 10384.2144
$ws.Range("N138").Value = -20664.2144

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H36").Value = 4905.2
$ws.Range("I36").Value = 1508.6666
$ws.Range("J36").Value = 10000
$ws.Range("K36").Value = 1508.6666
$ws.Range("L36").Value = 10000
$ws.Range("M36").Value = -1162.6666
$ws.Range("N36").Value = -10692

$ws.Range("H102").Value = 10420640
$ws.Range("I102").Value = 11907874
$ws.Range("K102").Value = 11907874
$ws.Range("M102").Value = -11906252

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 9101256
$ws.Range("I94").Value = 12989095
$ws.Range("K94").Value = 12989095
$ws.Range("M94").Value = -12988644

$ws.Range("H97").Value = 3542.2
$ws.Range("I97").Value = 3542.2
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 3542.2
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -2551.2
$ws.Range("N97").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 2399.8
$ws.Range("I62").Value = 1500
$ws.Range("K62").Value = 1500
$ws.Range("M62").Value = -876

$ws.Range("H65").Value = 2399.8
$ws.Range("I65").Value = 1500
$ws.Range("K65").Value = 7500
$ws.Range("M65").Value = -4380

$ws.Range("H122").Value = 4413.1665
$ws.Range("I122").Value = 4498.3335
$ws.Range("J122").Value = 4328
$ws.Range("K122").Value = 13495.0005
$ws.Range("L122").Value = 12984
$ws.Range("M122").Value = -11045.0005
$ws.Range("N122").Value = -17884

$ws.Range("H132").Value = 191000.1
$ws.Range("I132").Value = 168539.17
$ws.Range("K132").Value = 505617.51
$ws.Range("M132").Value = -503087.51

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 37080.07
$ws.Range("I5").Value = 765.2143
$ws.Range("J5").Value = 73394.92999999999
$ws.Range("K5").Value = 2295.6429
$ws.Range("L5").Value = 220184.79
$ws.Range("M5").Value = -2183.6429
$ws.Range("N5").Value = -220408.79

$ws.Range("H121").Value = 591.1111
$ws.Range("I121").Value = 122.14286
$ws.Range("K121").Value = 366.42858
$ws.Range("M121").Value = 943.57142

$ws.Range("H128").Value = 188855.28
$ws.Range("I128").Value = 188855.28
$ws.Range("K128").Value = 566565.84
$ws.Range("M128").Value = -561585.84

$ws.Range("H135").Value = 37080.07
$ws.Range("I135").Value = 765.2143
$ws.Range("J135").Value = 73394.92999999999
$ws.Range("K135").Value = 6886.928699999999
$ws.Range("L135").Value = 660554.3699999999
$ws.Range("M135").Value = -4351.928699999999
$ws.Range("N135").Value = -665624.3699999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4883804.5
$ws.Range("J80").Value = 8593.333000000001
$ws.Range("L80").Value = 8593.333000000001
$ws.Range("N80").Value = -10589.333

$ws.Range("H83").Value = 4883804.5
$ws.Range("J83").Value = 8593.333000000001
$ws.Range("L83").Value = 42966.665
$ws.Range("N83").Value = -52950.665

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 20000000
$ws.Range("I2").Value = 20000000
$ws.Range("K2").Value = 20000000
$ws.Range("M2").Value = -19999888

$ws.Range("H7").Value = 3409.8076
$ws.Range("I7").Value = 2020.6
$ws.Range("J7").Value = 5304.1816
$ws.Range("K7").Value = 2020.6
$ws.Range("L7").Value = 5304.1816
$ws.Range("M7").Value = -1908.6
$ws.Range("N7").Value = -5528.1816

$ws.Range("H74").Value = 47500
$ws.Range("J74").Value = 47000
$ws.Range("L74").Value = 47000
$ws.Range("N74").Value = -48996

$ws.Range("H77").Value = 47500
$ws.Range("J77").Value = 47000
$ws.Range("L77").Value = 141000
$ws.Range("N77").Value = -150984

$ws.Range("H122").Value = 4415.4194
$ws.Range("I122").Value = 2805.2632
$ws.Range("K122").Value = 8415.7896
$ws.Range("M122").Value = -5965.7896

$ws.Range("H126").Value = 3409.8076
$ws.Range("I126").Value = 2020.6
$ws.Range("J126").Value = 5304.1816
$ws.Range("K126").Value = 6061.799999999999
$ws.Range("L126").Value = 15912.5448
$ws.Range("M126").Value = -3591.799999999999
$ws.Range("N126").Value = -20852.5448

$ws.Range("H132").Value = 4431.109
$ws.Range("I132").Value = 3581.6216
$ws.Range("J132").Value = 6177.278
$ws.Range("K132").Value = 10744.8648
$ws.Range("L132").Value = 18531.834
$ws.Range("M132").Value = -8214.864799999999
$ws.Range("N132").Value = -23591.834

$ws.Range("H139").Value = 40715
$ws.Range("J139").Value = 40715
$ws.Range("L139").Value = 40715
$ws.Range("N139").Value = -50995

$ws.Range("H140").Value = 69714.25
$ws.Range("J140").Value = 69714.25
$ws.Range("L140").Value = 69714.25
$ws.Range("N140").Value = -80074.25

$ws.Range("H141").Value = 74001.92999999999
$ws.Range("J141").Value = 77859.21000000001
$ws.Range("L141").Value = 77859.21000000001
$ws.Range("N141").Value = -88219.21000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H11").Value = 3346661.8
$ws.Range("I11").Value = 3346661.8
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 3346661.8
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = -3346519.8
$ws.Range("N11").ClearContents()

$ws.Range("H62").Value = 6210.5254
$ws.Range("I62").Value = 2985.6072
$ws.Range("K62").Value = 2985.6072
$ws.Range("M62").Value = -2361.6072

$ws.Range("H65").Value = 6210.5254
$ws.Range("I65").Value = 2985.6072
$ws.Range("K65").Value = 14928.036
$ws.Range("M65").Value = -11808.036

$ws.Range("H75").Value = 15741.333
$ws.Range("I75").Value = 13614.5
$ws.Range("J75").Value = 19995
$ws.Range("K75").Value = 13614.5
$ws.Range("L75").Value = 19995
$ws.Range("M75").Value = -12678.5
$ws.Range("N75").Value = -21867

$ws.Range("H78").Value = 15741.333
$ws.Range("I78").Value = 13614.5
$ws.Range("J78").Value = 19995
$ws.Range("K78").Value = 40843.5
$ws.Range("L78").Value = 59985
$ws.Range("M78").Value = -36163.5
$ws.Range("N78").Value = -69345

$ws.Range("H113").Value = 704.6923
$ws.Range("I113").Value = 613.55554
$ws.Range("J113").Value = 909.75
$ws.Range("K113").Value = 1840.66662
$ws.Range("L113").Value = 2729.25
$ws.Range("M113").Value = 329.33338
$ws.Range("N113").Value = -7069.25

$ws.Range("H122").Value = 2325.2144
$ws.Range("I122").Value = 1099.6666
$ws.Range("J122").Value = 3244.375
$ws.Range("K122").Value = 3298.9998
$ws.Range("L122").Value = 9733.125
$ws.Range("M122").Value = -848.9998000000001
$ws.Range("N122").Value = -14633.125
